$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:G3").NumberFormat = "@"

$ws.Range("F1").Value = "latitude"
$ws.Range("G1").Value = "longitude"
$ws.Range("F2").Value = "0.735884"
$ws.Range("G2").Value = "34.5833"
$ws.Range("F3").Value = "0.055216"
$ws.Range("G3").Value = "34.2770"

$ws.Range("E1").Value = "Variety"
$ws.Range("E2").Value = "H84"
$ws.Range("E3").Value = "SC Saga"

$ws.Range("E1:G1").Font.Bold = $true

$ws.Columns.Item(7).ColumnWidth = 9.6

$ws.Range("E4").Select() | Out-Null
